$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ERD")
$ws2 = $wb.Worksheets.Item("Objects & Fields")

# ---------------------------------------------------------------------------
# 1) Transactions object: add a new "Transaction Date" / "Date" field as the
#    new last row of the group (row 47), and remove the thick bottom border
#    from the old last row ("Transaction_Status", row 46) since it is no
#    longer the last field of the group.
# ---------------------------------------------------------------------------
$ws2.Rows.Item(47).Insert()

# Row 47 is now a fresh blank row; give it the same border formatting that
# the group-ending row used to have (copied from row 46, which still carries
# the old "last row of group" formatting).
$ws2.Range("A46:D46").Copy()
$ws2.Range("A47:D47").PasteSpecial(-4122)
$ws2.Cells.Item(47, 2).Value = "Transaction Date"
$ws2.Cells.Item(47, 3).Value = "Date"
$ws2.Cells.Item(47, 4).ClearContents()

# Row 46 ("Transaction_Status") is no longer the last row in the group, so
# it loses its thick-bottom-border formatting and becomes a plain row.
$ws2.Range("A8:D8").Copy()
$ws2.Range("A46:D46").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Shipments object: add a new "Shipping_Date" / "Date" field as the new
#    last row of the group, and strip the thick bottom border from the old
#    last row ("Shipping_Charge").
#    After the insert above, "Shipping_Charge" now lives at row 54.
# ---------------------------------------------------------------------------
$ws2.Rows.Item(55).Insert()

$ws2.Range("A54:D54").Copy()
$ws2.Range("A55:D55").PasteSpecial(-4122)
$ws2.Cells.Item(55, 2).Value = "Shipping_Date"
$ws2.Cells.Item(55, 3).Value = "Date"
$ws2.Cells.Item(55, 4).ClearContents()

$ws2.Range("A3:D3").Copy()
$ws2.Range("A54:D54").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Box Themes object: the Box_ID auto-number format changes from
#    "B-{000000}" to "BI-{000000}". That row now lives at row 56.
# ---------------------------------------------------------------------------
$ws2.Cells.Item(56, 4).Value = "BI-{000000}"

# ---------------------------------------------------------------------------
# 4) Box Contents object: the Box_ID field length changes from 8 to 9.
#    That row now lives at row 60.
# ---------------------------------------------------------------------------
$ws2.Cells.Item(60, 4).Value = "Length: 9"

$ws2.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5) View/selection bookkeeping: update the selected cell on the
#    "Objects & Fields" sheet, then return focus to "ERD" so it remains the
#    active/selected tab (matching the saved file's original focus).
# ---------------------------------------------------------------------------
$ws2.Range("E55").Select()
$ws1.Activate()
$ws1.Range("U16").Select()
